$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting first for numeric-looking price strings so Excel keeps them as text
foreach ($addr in @("D5","D6","D7","D9","D10","D11","D12","D16","D17","D20","D21","D22","D23","D24","D25","D28","D29","D30","D31","D32","D33","D35","D37","D38","D39","D40","D41","D42","D45","D46","D47","D48","D49","D50","D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "44.706.46"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "2.238.10"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  +0.46%  "
$ws.Range("D5").Value = "305.67"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").Value = "95.01"
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("D7").Value = "0.571"
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").Value = "0.519"
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("D10").Value = "34.99"
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("D11").Value = "0.0802"
$ws.Range("E11").Value = "  -2.15%  "
$ws.Range("D12").Value = "7.18"
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "2.579.84"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").Value = "2.316.59"
$ws.Range("E15").Value = "  +2.99%  "
$ws.Range("D16").Value = "0.837"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "13.54"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").Value = "44.439.04"
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("D19").Value = "0.0₃0944"
$ws.Range("E19").Value = "  -3.26%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "6.25"
$ws.Range("E20").Value = "  -2.13%  "
$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").Value = "11.92"
$ws.Range("E21").Value = "  -2.41%  "
$ws.Range("D22").Value = "65.29"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").Value = "239.54"
$ws.Range("E23").Value = "  +0.97%  "
$ws.Range("D24").Value = "2.95"
$ws.Range("E24").Value = "  -0.48%  "
$ws.Range("D25").Value = "1.97"
$ws.Range("E25").Value = "  -1.66%  "
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("E27").Value = "  +4.71%  "
$ws.Range("D28").Value = "9.81"
$ws.Range("E28").Value = "  -1.60%  "
$ws.Range("D29").Value = "37.82"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "5.99"
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("D31").Value = "19.91"
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("D32").Value = "150.45"
$ws.Range("E32").Value = "  -1.22%  "
$ws.Range("D33").Value = "0.0792"
$ws.Range("E33").Value = "  -1.68%  "
$ws.Range("E34").Value = "  +0.94%  "
$ws.Range("D35").Value = "3.03"
$ws.Range("E35").Value = "  -8.91%  "
$ws.Range("E36").Value = "  -1.09%  "
$ws.Range("D37").Value = "0.107"
$ws.Range("E37").Value = "  -1.29%  "
$ws.Range("D38").Value = "1.84"
$ws.Range("E38").Value = "  +3.73%  "
$ws.Range("D39").Value = "15.11"
$ws.Range("E39").Value = "  +2.78%  "
$ws.Range("D40").Value = "3.37"
$ws.Range("E40").Value = "  -0.70%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.0301"
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "3.76"
$ws.Range("E42").Value = "  -3.05%  "
$ws.Range("E43").Value = "  +0.31%  "
$ws.Range("D44").Value = "1.832.15"
$ws.Range("E44").Value = "  +5.89%  "
$ws.Range("D45").Value = "1.73"
$ws.Range("E45").Value = "  +12.19%  "
$ws.Range("D46").Value = "79.54"
$ws.Range("E46").Value = "  -4.64%  "
$ws.Range("D47").Value = "0.189"
$ws.Range("E47").Value = "  -1.14%  "
$ws.Range("D48").Value = "98.53"
$ws.Range("E48").Value = "  -1.83%  "
$ws.Range("D49").Value = "4.88"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("D50").Value = "68.90"
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("D51").Value = "54.29"
$ws.Range("E51").Value = "  -0.98%  "

# Restore General number format now that the text values are stored
foreach ($addr in @("D5","D6","D7","D9","D10","D11","D12","D16","D17","D20","D21","D22","D23","D24","D25","D28","D29","D30","D31","D32","D33","D35","D37","D38","D39","D40","D41","D42","D45","D46","D47","D48","D49","D50","D51")) {
    $ws.Range($addr).NumberFormat = "General"
}
